$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: "79÷7=" -> "54÷6="
$t.Cell(1, 1).Range.Text = "54÷6="
# Row 1, Col 2: "33÷4=" -> "78÷9="
$t.Cell(1, 2).Range.Text = "78÷9="
# Row 1, Col 3: "74÷9=" -> "64÷7="
$t.Cell(1, 3).Range.Text = "64÷7="
# Row 1, Col 4: "34÷8=" -> "78÷8="
$t.Cell(1, 4).Range.Text = "78÷8="
# Row 1, Col 5: "32÷9=" -> "13÷4="
$t.Cell(1, 5).Range.Text = "13÷4="
# Row 5, Col 1: "34÷7=" -> "22÷8="
$t.Cell(5, 1).Range.Text = "22÷8="
# Row 5, Col 2: "58÷2=" -> "21÷2="
$t.Cell(5, 2).Range.Text = "21÷2="
# Row 5, Col 3: "60÷8=" -> "25÷8="
$t.Cell(5, 3).Range.Text = "25÷8="
# Row 5, Col 4: "63÷5=" -> "12÷8="
$t.Cell(5, 4).Range.Text = "12÷8="
# Row 5, Col 5: "22÷7=" -> "86÷3="
$t.Cell(5, 5).Range.Text = "86÷3="
# Row 9, Col 1: "73÷7=" -> "17÷2="
$t.Cell(9, 1).Range.Text = "17÷2="
# Row 9, Col 2: "98÷7=" -> "83÷8="
$t.Cell(9, 2).Range.Text = "83÷8="
# Row 9, Col 3: "79÷3=" -> "26÷8="
$t.Cell(9, 3).Range.Text = "26÷8="
# Row 9, Col 4: "65÷2=" -> "19÷4="
$t.Cell(9, 4).Range.Text = "19÷4="
# Row 9, Col 5: "46÷4=" -> "89÷2="
$t.Cell(9, 5).Range.Text = "89÷2="
# Row 13, Col 1: "39÷2=" -> "95÷5="
$t.Cell(13, 1).Range.Text = "95÷5="
# Row 13, Col 2: "95÷3=" -> "34÷2="
$t.Cell(13, 2).Range.Text = "34÷2="
# Row 13, Col 3: "11÷9=" -> "22÷7="
$t.Cell(13, 3).Range.Text = "22÷7="
# Row 13, Col 4: "35÷8=" -> "38÷2="
$t.Cell(13, 4).Range.Text = "38÷2="
# Row 13, Col 5: "42÷6=" -> "29÷3="
$t.Cell(13, 5).Range.Text = "29÷3="
# Row 17, Col 1: "11÷4=" -> "91÷7="
$t.Cell(17, 1).Range.Text = "91÷7="
# Row 17, Col 2: "30÷9=" -> "19÷6="
$t.Cell(17, 2).Range.Text = "19÷6="
# Row 17, Col 3: "26÷8=" -> "40÷6="
$t.Cell(17, 3).Range.Text = "40÷6="
# Row 17, Col 4: "80÷2=" -> "93÷2="
$t.Cell(17, 4).Range.Text = "93÷2="
# Row 17, Col 5: "23÷6=" -> "94÷7="
$t.Cell(17, 5).Range.Text = "94÷7="
